$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.596.73'
$ws.Range('E2').Value = '  +14.23%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.823.48'
$ws.Range('E3').Value = '  +8.78%  '

# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.994'
$ws.Range('E4').Value = '  -0.52%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.40'
$ws.Range('E5').Value = '  +6.75%  '

# Row 6: XRP
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.553'
$ws.Range('E6').Value = '  +5.73%  '

# Row 7: USDC
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.20%  '

# Row 8: Solana
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.03'
$ws.Range('E8').Value = '  +7.46%  '

# Row 9: OKB
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.30'
$ws.Range('E9').Value = '  +5.25%  '

# Row 10: Cardano
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.286'
$ws.Range('E10').Value = '  +7.84%  '

# Row 11: Dogecoin
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0685'
$ws.Range('E11').Value = '  +10.67%  '

# Row 12: TRON
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0929'
$ws.Range('E12').Value = '  +2.50%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.060.73'
$ws.Range('E13').Value = '  +7.50%  '

# Row 14: WrappedEther
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.809.62'
$ws.Range('E14').Value = '  +8.03%  '

# Row 15: Polygon
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.648'
$ws.Range('E15').Value = '  +4.52%  '

# Row 16: WrappedBTC
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.498.98'
$ws.Range('E16').Value = '  +13.83%  '

# Row 17: Chainlink
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '10.37'
$ws.Range('E17').Value = '  -3.15%  '

# Row 18: Polkadot
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.37'
$ws.Range('E18').Value = '  +9.02%  '

# Row 19: Litecoin
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.66'
$ws.Range('E19').Value = '  +8.98%  '

# Row 20: BitcoinCash
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '265.58'
$ws.Range('E20').Value = '  +7.28%  '

# Row 21: ShibaInu
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0766'
$ws.Range('E21').Value = '  +6.52%  '

# Row 22: Dai
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('E22').Value = '  -0.21%  '

# Row 23: Uniswap
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.44'
$ws.Range('E23').Value = '  +2.91%  '

# Row 24: Avalanche
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.54'
$ws.Range('E24').Value = '  +4.41%  '

# Row 25: Toncoin
$ws.Range('E25').Value = '  -1.97%  '

# Row 26: Monero
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.36'
$ws.Range('E26').Value = '  +2.18%  '

# Row 27: EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.06'
$ws.Range('E27').Value = '  +7.47%  '

# Row 28: Stellar
$ws.Range('E28').Value = '  +5.46%  '

# Row 29: Cosmos
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.19'
$ws.Range('E29').Value = '  +6.29%  '

# Row 30: BinanceUSD
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.994'
$ws.Range('E30').Value = '  -0.47%  '

# Row 31: Filecoin
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.90'
$ws.Range('E31').Value = '  +12.15%  '

# Row 32: Hedera
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0519'
$ws.Range('E32').Value = '  +3.76%  '

# Row 33: PancakeSwap
$ws.Range('E33').Value = '  +6.73%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.62'
$ws.Range('E34').Value = '  +9.62%  '

# Row 35: Maker
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.586.14'
$ws.Range('E35').Value = '  +7.06%  '

# Row 36: LidoDAOToken
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.86'
$ws.Range('E36').Value = '  +7.86%  '

# Row 37: Aave
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '89.87'
$ws.Range('E37').Value = '  +12.96%  '

# Row 38: TrustWalletToken
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').Value = '  +3.25%  '

# Row 39: ImmutableX
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.632'
$ws.Range('E39').Value = '  +7.34%  '

# Row 40: VeChain
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0189'
$ws.Range('E40').Value = '  +5.79%  '

# Row 41: MXToken
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('E41').Value = '  +6.35%  '

# Row 42: ARBITRUM
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.932'
$ws.Range('E42').Value = '  +8.41%  '

# Row 43: HuobiToken
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.36'
$ws.Range('E43').Value = '  +2.68%  '

# Row 44: RenderToken
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.16'
$ws.Range('E44').Value = '  +6.90%  '

# Row 45: Kaspa
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0522'
$ws.Range('E45').Value = '  +3.24%  '

# Row 46: WEMIXToken
$ws.Range('E46').Value = '  +2.83%  '

# Row 47: RocketPoolETH
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.954.17'
$ws.Range('E47').Value = '  +7.69%  '

# Row 48: BitcoinSV
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.28'
$ws.Range('E48').Value = '  +3.68%  '

# Row 49: FraxShare
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.80'
$ws.Range('E49').Value = '  +6.56%  '

# Row 50: PaxDollar
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  +0.06%  '

# Row 51: InjectiveProtocol
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.56'
$ws.Range('E51').Value = '  +23.84%  '
